$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132 (Item ID 44049)
$ws.Range("H132").Value = 4724.7446
$ws.Range("I132").Value = 3204.8975
$ws.Range("J132").Value = 12134
$ws.Range("K132").Value = 9614.692500000001
$ws.Range("L132").Value = 36402
$ws.Range("M132").Value = -7084.692500000001
$ws.Range("N132").Value = -41462

# Row 137 (Item ID 44013)
$ws.Range("H137").Value = 3484.3
$ws.Range("I137").Value = 3098.7693
$ws.Range("J137").Value = 4200.2856
$ws.Range("K137").Value = 9296.3079
$ws.Range("L137").Value = 12600.8568
$ws.Range("M137").Value = -6746.3079
$ws.Range("N137").Value = -17700.8568

# Row 138 (Item ID 44169)
$ws.Range("H138").Value = 1789.8868
$ws.Range("I138").Value = 1009.9487
$ws.Range("J138").Value = 3962.5715
$ws.Range("K138").Value = 3029.8461
$ws.Range("L138").Value = 11887.7145
$ws.Range("M138").Value = 2110.1539
$ws.Range("N138").Value = -22167.7145

# Row 141 (Item ID 44161)
$ws.Range("H141").Value = 1299.2759
$ws.Range("I141").Value = 1395.8846
$ws.Range("J141").Value = 462
$ws.Range("K141").Value = 4187.6538
$ws.Range("L141").Value = 1386
$ws.Range("M141").Value = 992.3462
$ws.Range("N141").Value = -11746

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Item ID 27713)
$ws.Range("H2").Value = 1765.069
$ws.Range("I2").Value = 1507.2307
$ws.Range("K2").Value = 1507.2307
$ws.Range("M2").Value = -1394.2307

# Row 44 (Item ID 3861)
$ws.Range("H44").Value = 49499.5
$ws.Range("J44").Value = 49499.5
$ws.Range("L44").Value = 49499.5
$ws.Range("N44").Value = -50475.5

# Row 61 (Item ID 43999)
$ws.Range("H61").Value = 2724.647
$ws.Range("I61").Value = 2687.9333
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2687.9333
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2475.9333
$ws.Range("N61").Value = -3424

# Row 101 (Item ID 18518)
$ws.Range("H101").Value = 49798.5
$ws.Range("J101").Value = 49798.5
$ws.Range("L101").Value = 49798.5
$ws.Range("N101").Value = -56288.5

# Row 116 (Item ID 27713)
$ws.Range("H116").Value = 1765.069
$ws.Range("I116").Value = 1507.2307
$ws.Range("K116").Value = 1507.2307
$ws.Range("M116").Value = 786.7692999999999

# Row 136 (Item ID 43999)
$ws.Range("H136").Value = 2724.647
$ws.Range("I136").Value = 2687.9333
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 8063.7999
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -5513.7999
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Item ID 27713)
$ws.Range("H3").Value = 1765.069
$ws.Range("I3").Value = 1507.2307
$ws.Range("K3").Value = 1507.2307
$ws.Range("M3").Value = -1393.2307

# Row 134 (Item ID 43998)
$ws.Range("H134").Value = 2786.48
$ws.Range("I134").Value = 2568
$ws.Range("K134").Value = 7704
$ws.Range("M134").Value = -5169

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Item ID 44023)
$ws.Range("H31").Value = 3168.3513
$ws.Range("I31").Value = 2805.3076
$ws.Range("J31").Value = 4026.4546
$ws.Range("K31").Value = 2805.3076
$ws.Range("L31").Value = 4026.4546
$ws.Range("M31").Value = -2510.3076
$ws.Range("N31").Value = -4616.4546

# Row 34 (Item ID 44023)
$ws.Range("H34").Value = 3168.3513
$ws.Range("I34").Value = 2805.3076
$ws.Range("J34").Value = 4026.4546
$ws.Range("K34").Value = 2805.3076
$ws.Range("L34").Value = 4026.4546
$ws.Range("M34").Value = -2603.3076
$ws.Range("N34").Value = -4430.4546

# Row 43 (Item ID 18504)
$ws.Range("H43").Value = 37413.75
$ws.Range("J43").Value = 37413.75
$ws.Range("L43").Value = 37413.75
$ws.Range("N43").Value = -37781.75

# Row 101 (Item ID 18504)
$ws.Range("H101").Value = 37413.75
$ws.Range("J101").Value = 37413.75
$ws.Range("L101").Value = 37413.75
$ws.Range("N101").Value = -43903.75

# Row 102 (Item ID 19738)
$ws.Range("H102").Value = 79999
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 79999
$ws.Range("K102").Value = 0
$ws.Range("L102").ClearContents()
$ws.Range("M102").Value = 79999
$ws.Range("N102").Value = -84867

# Row 103 (Item ID 19558)
$ws.Range("H103").Value = 48704
$ws.Range("I103").Value = 43380.25
$ws.Range("K103").Value = 43380.25
$ws.Range("M103").Value = -42208.25

# Row 104 (Item ID 19749)
$ws.Range("H104").Value = 97285
$ws.Range("J104").Value = 97285
$ws.Range("L104").Value = 97285
$ws.Range("N104").Value = -102527

# Row 105 (Item ID 19928)
$ws.Range("H105").Value = 40839
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# Row 106 (Item ID 18661)
$ws.Range("H106").Value = 78900
$ws.Range("J106").Value = 78900
$ws.Range("L106").Value = 78900
$ws.Range("N106").Value = -81424

# Row 134 (Item ID 44020)
$ws.Range("H134").Value = 19675.264
$ws.Range("I134").Value = 24748.363
$ws.Range("J134").Value = 2504.7693
$ws.Range("K134").Value = 74245.08900000001
$ws.Range("L134").Value = 7514.3079
$ws.Range("M134").Value = -71710.08900000001
$ws.Range("N134").Value = -12584.3079

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (Item ID 43974)
$ws.Range("H5").Value = 686.93335
$ws.Range("I5").Value = 741.1
$ws.Range("J5").Value = 578.6
$ws.Range("K5").Value = 2223.3
$ws.Range("L5").Value = 1735.8
$ws.Range("M5").Value = -2111.3
$ws.Range("N5").Value = -1959.8

# Row 12 (Item ID 4854)
$ws.Range("H12").Value = 311.16666
$ws.Range("I12").Value = 266.66666
$ws.Range("J12").Value = 355.66666
$ws.Range("K12").Value = 799.9999799999999
$ws.Range("L12").Value = 1066.99998
$ws.Range("M12").Value = -626.9999799999999
$ws.Range("N12").Value = -1412.99998

# Row 39 (Item ID 4712)
$ws.Range("H39").Value = 8696.454
$ws.Range("J39").Value = 8696.454
$ws.Range("L39").Value = 26089.362
$ws.Range("N39").Value = -26677.362

# Row 60 (Item ID 4750)
$ws.Range("H60").Value = 206.33333
$ws.Range("I60").Value = 219.5
$ws.Range("J60").Value = 180
$ws.Range("K60").Value = 658.5
$ws.Range("L60").Value = 540
$ws.Range("M60").Value = -407.5
$ws.Range("N60").Value = -1042

# Row 75 (Item ID 12863)
$ws.Range("H75").Value = 725.6667
$ws.Range("I75").Value = 700
$ws.Range("J75").Value = 738.5
$ws.Range("K75").Value = 2100
$ws.Range("L75").Value = 2215.5
$ws.Range("M75").Value = -1102
$ws.Range("N75").Value = -4211.5

# Row 78 (Item ID 12863)
$ws.Range("H78").Value = 725.6667
$ws.Range("I78").Value = 700
$ws.Range("J78").Value = 738.5
$ws.Range("K78").Value = 6300
$ws.Range("L78").Value = 6646.5
$ws.Range("M78").Value = -1308
$ws.Range("N78").Value = -16630.5

# Row 109 (Item ID 27854)
$ws.Range("H109").Value = 2129.7
$ws.Range("I109").Value = 2033
$ws.Range("K109").Value = 6099
$ws.Range("M109").Value = -5059

# Row 135 (Item ID 43974)
$ws.Range("H135").Value = 686.93335
$ws.Range("I135").Value = 741.1
$ws.Range("J135").Value = 578.6
$ws.Range("K135").Value = 6669.900000000001
$ws.Range("L135").Value = 5207.400000000001
$ws.Range("M135").Value = -4134.900000000001
$ws.Range("N135").Value = -10277.4

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Item ID 14146)
$ws.Range("H70").Value = 8387.16
$ws.Range("I70").Value = 9091.714
$ws.Range("J70").Value = 7490.4546
$ws.Range("K70").Value = 9091.714
$ws.Range("L70").Value = 7490.4546
$ws.Range("M70").Value = -8821.714
$ws.Range("N70").Value = -8030.4546

# Row 73 (Item ID 14146)
$ws.Range("H73").Value = 8387.16
$ws.Range("I73").Value = 9091.714
$ws.Range("J73").Value = 7490.4546
$ws.Range("K73").Value = 9091.714
$ws.Range("L73").Value = 7490.4546
$ws.Range("M73").Value = -8155.714
$ws.Range("N73").Value = -9362.454600000001

# Row 97 (Item ID 19940)
$ws.Range("H97").Value = 3895.6667
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 3895.6667
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("M97").Value = 3895.6667
$ws.Range("N97").Value = -4887.6667

# Row 102 (Item ID 36169)
$ws.Range("H102").Value = 4504.3335
$ws.Range("I102").Value = 624.76
$ws.Range("J102").Value = 52999
$ws.Range("K102").Value = 624.76
$ws.Range("L102").Value = 52999
$ws.Range("M102").Value = 997.24
$ws.Range("N102").Value = -56243

# Row 132 (Item ID 44008)
$ws.Range("H132").Value = 26387.666
$ws.Range("I132").Value = 30212.139
$ws.Range("J132").Value = 3440.8333
$ws.Range("K132").Value = 90636.417
$ws.Range("L132").Value = 10322.4999
$ws.Range("M132").Value = -88106.417
$ws.Range("N132").Value = -15382.4999

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Item ID 5277)
$ws.Range("H22").Value = 1240
$ws.Range("I22").Value = 1240
$ws.Range("K22").Value = 1240
$ws.Range("M22").Value = -945

# Row 27 (Item ID 5277)
$ws.Range("H27").Value = 1240
$ws.Range("I27").Value = 1240
$ws.Range("K27").Value = 1240
$ws.Range("M27").Value = -1133

# Row 46 (Item ID 5282)
$ws.Range("H46").Value = 6639.1763
$ws.Range("I46").Value = 14121
$ws.Range("J46").Value = 1401.9
$ws.Range("K46").Value = 14121
$ws.Range("L46").Value = 1401.9
$ws.Range("M46").Value = -13933
$ws.Range("N46").Value = -1777.9

# Row 55 (Item ID 5284)
$ws.Range("H55").Value = 2925
$ws.Range("I55").Value = 3066.6667
$ws.Range("J55").Value = 2500
$ws.Range("K55").Value = 3066.6667
$ws.Range("L55").Value = 2500
$ws.Range("M55").Value = -2893.6667
$ws.Range("N55").Value = -2846

# Row 103 (Item ID 18526)
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("N103").Value = 0

# Row 118 (Item ID 26146)
$ws.Range("H118").Value = 69000
$ws.Range("J118").Value = 69000
$ws.Range("L118").Value = 69000
$ws.Range("N118").Value = -72314

# Row 122 (Item ID 36247)
$ws.Range("H122").Value = 193421.61
$ws.Range("I122").Value = 5030.5
$ws.Range("K122").Value = 15091.5
$ws.Range("M122").Value = -12641.5

# Row 132 (Item ID 44058)
$ws.Range("H132").Value = 26561.824
$ws.Range("I132").Value = 29339.178
$ws.Range("K132").Value = 88017.534
$ws.Range("M132").Value = -85487.534

$ws = $wb.Worksheets.Item("WVR")
# Row 15 (Item ID 2670)
$ws.Range("H15").Value = 91650
$ws.Range("I15").Value = 91650
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 91650
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -91362
